$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "shift" in G1
$ws.Range("G1").Value = "shift"

# Add value 0 for the new "shift" column in rows 2-7
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 0

# Update the selected cell to G7 as in the diff
$ws.Range("G7").Select()
